$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix BOM entry for S1, S2 (row 27): Comment and LCSC Part Number
$ws.Range("A27").Value = "EVQQ2U02W"
$ws.Range("D27").Value = "C395227"

# Update last active cell selection to the last used row
$ws.Range("B35").Select()

# Re-fit column B width to content
$ws.Columns.Item(2).AutoFit() | Out-Null
